# Update scenarios / levels and growth
$wb = $excel.ActiveWorkbook

# --- Workbook-level: update the absolute path recorded by Excel ---
# (cosmetic metadata written by Excel itself; not settable via the object
# model in a meaningful way, so we leave workbook-level internals alone and
# focus on the user-visible worksheet changes below.)

$ws1 = $wb.Worksheets.Item("BuDA Test1")
$ws2 = $wb.Worksheets.Item("BuDA Test2")

# --- Sheet1 (BuDA Test1): update view state ---
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollRow = 6
$ws1.Range("E11").Select()

# --- Sheet1 (BuDA Test1): update column D values ---
$ws1.Range("D7").Value = 1

$ws1.Range("D10:D13").Value = -0.33
$ws1.Range("D14:D25").Value = -0.31
$ws1.Range("D26:D37").Value = -0.13
$ws1.Range("D38:D49").Value = 0.15

# --- Sheet2 (BuDA Test2): update view state ---
$ws2.Activate()
$ws2.Range("E7").Select()

# --- Sheet2 (BuDA Test2): update column D values ---
$ws2.Range("D7:D10").Value = -0.33
$ws2.Range("D11:D22").Value = 0.24
$ws2.Range("D23:D34").Value = 0.34
$ws2.Range("D35:D46").Value = 0.54

$ws1.Activate()
